# Replace the city placeholders ("Saint Malo" / "Saint-Malo") with the
# {{organisme_ville}} template variable and, while doing so, collapse each
# of the two affected paragraphs (which were split across several runs
# with identical/near-identical formatting) down to a single plain run
# with no run-level formatting at all — mirroring how the author's
# templating pass rewrote these paragraphs wholesale.

$d = $word.ActiveDocument

# --- ARTICLE 7: "Si une contestation ... le litige." -----------------
$r1 = $d.Content
$r1.Find.Execute("Si une contestation ou un différend", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1.Expand(4) | Out-Null          # wdParagraph -> grab the whole paragraph (incl. mark)
$target1 = $d.Range($r1.Start, $r1.End - 1)   # exclude the trailing paragraph mark
$target1.Delete()
$target1.InsertAfter("Si une contestation ou un différend n’a pu être réglé à l’amiable, le tribunal de commerce de {{organisme_ville}} sera seul compétent pour régler le litige.")

# --- Signature block: "Fait en double exemplaire, ..." ----------------
$r2 = $d.Content
$r2.Find.Execute("Fait en double exemplaire", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Expand(4) | Out-Null
$target2 = $d.Range($r2.Start, $r2.End - 1)
$target2.Delete()
$target2.InsertAfter("Fait en double exemplaire, à {{organisme_ville}}, le 5 septembre 2022,")
